# "Add files via upload" - daily snapshot update.
#
# 1) Sheet3!B20:B36 holds the latest per-product lookup values that feed
#    Sheet3!C2:C18 (IFERROR/VLOOKUP) and Sheet1's CB/CC VLOOKUP columns.
#    Refresh those raw values for today.
# 2) Sheet1 gets a new trailing snapshot column (CN) for "22-nov", holding
#    a copy of the (now refreshed) CB/CC value for each product row.
# 3) The active selection moves to the new column.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- 1) Refresh Sheet3's raw lookup table (A20:B36) ------------------------
$newValues = @{
    20 = 7.3830031133737029   # 3D QUESO 92GX27
    21 = 2.3321091203297164   # CHEETOS QUESO 85GX24X1
    22 = 7.6975611876011918   # DORITOS QUESO 129GX19
    23 = 8.942351403102295    # DORITOS QUESO 70X40G
    24 = 6.990309759752801    # DORITOS QUESO 77GX26
    26 = 14.566260958152959   # LAYS CLASICAS 145GRX18
    27 = 10.026785466472022   # LAYS CLASICAS 249GRX14
    28 = 5.2875743553630024   # LAYS CLASICAS 40GX68
    29 = 0.27660253877070262  # LAYS CLASICAS 94GRX25
    31 = 11.938850000099501   # LAYS ONDAS FH 70GX28
    32 = 12.042802705067867   # LAYS QSO Y CEBOLLA 34GX72
    33 = 10.458726400575104   # PEHUAMAR ACANALADA 520GX9
    34 = 12.550918272060255   # PEHUAMAR MAICITOS 285GX10
    35 = 7.1956910702905352   # PEHUAMAR PAPA LISA 520GX9
    36 = 34.320681219573011   # QUAKER AVENA INSTANT FORTIF 18X280G
}

foreach ($r in $newValues.Keys) {
    $ws3.Range("B" + $r).Value = $newValues[$r]
}

# --- 2) Add the new "22-nov" snapshot column on Sheet1 ---------------------
$ws1.Range("CN1").Value = "22-nov"
$ws1.Range("CN1").NumberFormat = $ws1.Range("CM1").NumberFormat()

for ($r = 2; $r -le 18; $r++) {
    $src = $ws1.Range("CB" + $r).Value()
    $dst = $ws1.Range("CN" + $r)
    $dst.Value = $src
    $dst.NumberFormat = $ws1.Range("CM" + $r).NumberFormat()
}

# --- 3) Move the selection onto the freshly added column -------------------
$ws1.Activate()
$ws1.Range("CN4").Select() | Out-Null
